# Trade #23 closed at 2026-02-17 23:57:23 - unknown UNKNOWN +0.000%
#
# Updates the "Summary" and "Strategy Status" roll-up sheets for the newly
# closed MarketMaking trade, and appends the trade's row to both the
# "All Trades" and "MarketMaking" ledgers.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Helper: write a text value to a cell while preventing Excel's COM layer
# from auto-coercing date-/time-looking strings (e.g. "2026-02-17",
# "23:57:16") into date/time serial numbers. We briefly force a Text
# number format for the assignment, then clear the explicit formatting
# so the cell is left with the default (unstyled) look, matching the
# other rows in the sheet.
# ---------------------------------------------------------------------
function Set-TextValue($cell, $text) {
    $cell.NumberFormat = "@"
    $cell.Value = $text
    $cell.ClearFormats()
}

# ---------------------------------------------------------------------
# Summary sheet
# ---------------------------------------------------------------------
$summary = $wb.Worksheets.Item("Summary")
$summary.Range("B3").Value = 1500.83   # Current Capital
$summary.Range("B4").Value = 0.83      # Total P&L $
$summary.Range("B5").Value = 0.72      # Total P&L %
$summary.Range("B6").Value = 23        # Total Trades
$summary.Range("B7").Value = 13        # Winning Trades
$summary.Range("B9").Value = 56.52     # Win Rate %

# ---------------------------------------------------------------------
# Strategy Status sheet - MarketMaking row (row 6)
# ---------------------------------------------------------------------
$status = $wb.Worksheets.Item("Strategy Status")
$status.Range("C6").Value = 100.83     # Capital
$status.Range("D6").Value = 23         # Trades
$status.Range("E6").Value = 0.83       # P&L $
$status.Range("F6").Value = 0.83       # P&L %
$status.Range("G6").Value = 56.52      # Win Rate %

# ---------------------------------------------------------------------
# Append the new closed trade (row 24) to both the "All Trades" and
# "MarketMaking" ledgers - they mirror each other since MarketMaking is
# the only currently active strategy.
# ---------------------------------------------------------------------
$sheetNames = @("All Trades", "MarketMaking")

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)

    $ws.Cells.Item(24, 1).Value = 23                 # Trade #
    Set-TextValue $ws.Cells.Item(24, 2) "2026-02-17" # Date
    Set-TextValue $ws.Cells.Item(24, 3) "23:57:16"   # Time
    Set-TextValue $ws.Cells.Item(24, 4) "MarketMaking" # Strategy
    Set-TextValue $ws.Cells.Item(24, 5) "UP"         # Side
    $ws.Cells.Item(24, 6).Value = 0.5600000000000001 # Entry Price
    $ws.Cells.Item(24, 7).Value = 0.64                # Exit Price
    Set-TextValue $ws.Cells.Item(24, 8) "CLOSED"     # Status
    $ws.Cells.Item(24, 9).Value = 14.2857             # P&L %
    $ws.Cells.Item(24, 10).Value = 0.08               # P&L $
    $ws.Cells.Item(24, 11).Value = 100.83             # Capital After
    $ws.Cells.Item(24, 12).Value = 0                  # Entry Slippage (bps)
    $ws.Cells.Item(24, 13).Value = 0                  # Exit Slippage (bps)
    $ws.Cells.Item(24, 14).Value = 0.6                # Confidence
    Set-TextValue $ws.Cells.Item(24, 15) "Normal spread capture: 19600 bps" # Entry Reason
    Set-TextValue $ws.Cells.Item(24, 16) "early_exit" # Exit Reason
    $ws.Cells.Item(24, 17).Value = 0.13               # Duration (min)
}
